# SESSION資料結構.xlsx — add new rows describing select_course / login flow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: section header "抓select_course的資料(course)" (same yellow style as A3/A8/A15) ---
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(22,1).PasteSpecial(-4122)
$ws.Cells.Item(22,1).Value = "抓select_course的資料(course)"

# --- Row 26: 登入(login) header | 判斷是學生還是老師 note ---
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(26,1).PasteSpecial(-4122)
$ws.Cells.Item(26,1).Value = "登入(login)"
$ws.Cells.Item(13,4).Copy()
$ws.Cells.Item(26,4).PasteSpecial(-4122)
$ws.Cells.Item(26,4).Value = "判斷是學生還是老師"

# --- Row 23: 各種unit_id | 數字0~N | 存unit_course第N個學生的id ---
$ws.Cells.Item(23,2).Value = "各種unit_id"
$ws.Cells.Item(23,3).Value = "數字0~N"
$ws.Cells.Item(13,4).Copy()
$ws.Cells.Item(23,4).PasteSpecial(-4122)
$ws.Cells.Item(23,4).Value = "存unit_course第N個學生的id"

# --- Row 24: 各種unit_id | 第N個學生的stu_id | 這個select_course的id ---
$ws.Cells.Item(24,2).Value = "各種unit_id"
$ws.Cells.Item(24,3).Value = "第N個學生的stu_id"
$ws.Cells.Item(13,4).Copy()
$ws.Cells.Item(24,4).PasteSpecial(-4122)
$ws.Cells.Item(24,4).Value = "這個select_course的id"

$excel.CutCopyMode = 0

# --- Column A got wider to fit the new, longer labels ---
$ws.Columns.Item(1).ColumnWidth = 28.7

# --- Selection moved to where the user was last editing ---
$ws.Cells.Item(27,2).Select()
